$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 71
$ws1.Range("F5").Value = 362
$ws1.Range("F6").Value = 5531
$ws1.Range("F8").Value = 6439
$ws1.Range("F9").Value = 643
$ws1.Range("F10").Value = 12
$ws1.Range("F11").Value = 1403
$ws1.Range("F12").Value = 40
$ws1.Range("F13").Value = 110

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 71
$ws4.Range("F6").Value = 362
$ws4.Range("F7").Value = 5531
$ws4.Range("F9").Value = 6439
$ws4.Range("F10").Value = 643
$ws4.Range("F11").Value = 12
$ws4.Range("F12").Value = 1403
$ws4.Range("F13").Value = 40
$ws4.Range("F14").Value = 110
